$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a scratch cell in a column that has no column-level style to mint a
# "clean" font/style pair (Excel drops the redundant charset attribute when
# it re-declares the default font after an explicit font assignment).
$temp = $ws.Cells.Item(1, 10)
$temp.Value = "x"
$temp.Font.Name = "Calibri"
$temp.Copy()

$row = 64
$ws.Cells.Item($row, 1).Value = "Albert Rapp"
$ws.Cells.Item($row, 2).Value = "Creating beautiful tables in R with {gt}"
$ws.Cells.Item($row, 3).Value = "https://gt.albert-rapp.de/"

$ws.Range("A64:C64").PasteSpecial(-4122)  # xlPasteFormats

$temp.Clear()

$ws.Range("A64:C64").Select()
